$wb = $excel.ActiveWorkbook

# --- DIVE_SITE_METADATA: rename "L#" -> "Locality Code", "Site#" -> "Site Code", "Functions:" -> "Functions"
$ws1 = $wb.Worksheets.Item("DIVE_SITE_METADATA")
$ws1.Range("D1").Value = "Locality Code"
$ws1.Range("F1").Value = "Site Code"
$ws1.Range("S1").Value = "Functions"
$ws1.Columns.Item(4).ColumnWidth = 12.74
$ws1.Columns.Item(6).ColumnWidth = 8.89
$ws1.Range("S2").Select()

# --- BENTHIC_TAXAS: rename "Indicators:" -> "Indicators"
$ws3 = $wb.Worksheets.Item("BENTHIC_TAXAS")
$ws3.Range("E1").Value = "Indicators"
$ws3.Range("E2").Select()

# --- MOTILE_DB: insert new "Surveyed area" column before "Size"
$ws4 = $wb.Worksheets.Item("MOTILE_DB")
$ws4.Columns.Item(5).Insert()
$ws4.Range("E1").Value = "Surveyed area"
$ws4.Columns.Item(5).ColumnWidth = 13.31
$ws4.Range("H3").Select()

# --- MOTILE_TAXAS: rename "Indicators:" -> "Indicators"
$ws5 = $wb.Worksheets.Item("MOTILE_TAXAS")
$ws5.Range("E1").Value = "Indicators"
$ws5.Range("E2").Select()

$ws1.Select()
